# Auto-generated Excel COM-interop script
# Applies the reordering of same-date match rows and appends 3 new match rows
# as described by the target diff (commit: 'Atualizado por script em 12-11-2023 14:45').

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: re-write columns F:V for rows whose match data (same kickoff date)
# was re-ordered by the upstream scraper. Columns A:E (Indice/pais/torneio/
# temporada/data_partida) are untouched for these rows.
# ---------------------------------------------------------------------------

# Row 2
$ws.Cells.Item(2, 6).Value = 'Smederevo'
$ws.Cells.Item(2, 7).Value = 2
$ws.Cells.Item(2, 8).Value = 'FK Indjija'
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 2.74
$ws.Cells.Item(2, 11).Value = '05/08/2023 13:12'
$ws.Cells.Item(2, 12).Value = 2.74
$ws.Cells.Item(2, 13).Value = '05/08/2023 13:12'
$ws.Cells.Item(2, 14).Value = 2.75
$ws.Cells.Item(2, 15).Value = '05/08/2023 13:12'
$ws.Cells.Item(2, 16).Value = 2.76
$ws.Cells.Item(2, 17).Value = '05/08/2023 15:33'
$ws.Cells.Item(2, 18).Value = 2.56
$ws.Cells.Item(2, 19).Value = '05/08/2023 13:12'
$ws.Cells.Item(2, 20).Value = 2.56
$ws.Cells.Item(2, 21).Value = '05/08/2023 13:12'
$ws.Cells.Item(2, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/smederevo-indjija/rR4gggd8/'

# Row 3
$ws.Cells.Item(3, 6).Value = 'Vrsac'
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = 'Kolubara'
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 2.71
$ws.Cells.Item(3, 11).Value = '05/08/2023 13:12'
$ws.Cells.Item(3, 12).Value = 2.71
$ws.Cells.Item(3, 13).Value = '05/08/2023 13:12'
$ws.Cells.Item(3, 14).Value = 2.81
$ws.Cells.Item(3, 15).Value = '05/08/2023 13:12'
$ws.Cells.Item(3, 16).Value = 2.83
$ws.Cells.Item(3, 17).Value = '05/08/2023 15:33'
$ws.Cells.Item(3, 18).Value = 2.53
$ws.Cells.Item(3, 19).Value = '05/08/2023 13:12'
$ws.Cells.Item(3, 20).Value = 2.53
$ws.Cells.Item(3, 21).Value = '05/08/2023 13:12'
$ws.Cells.Item(3, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/vrsac-kolubara/rVQRcbII/'

# Row 4
$ws.Cells.Item(4, 6).Value = 'Jedinstvo U.'
$ws.Cells.Item(4, 7).Value = 4
$ws.Cells.Item(4, 8).Value = 'RFK Novi Sad'
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 1.76
$ws.Cells.Item(4, 11).Value = '04/08/2023 05:42'
$ws.Cells.Item(4, 12).Value = 1.44
$ws.Cells.Item(4, 13).Value = '05/08/2023 17:25'
$ws.Cells.Item(4, 14).Value = 3.15
$ws.Cells.Item(4, 15).Value = '04/08/2023 05:42'
$ws.Cells.Item(4, 16).Value = 4.59
$ws.Cells.Item(4, 17).Value = '05/08/2023 17:25'
$ws.Cells.Item(4, 18).Value = 3.94
$ws.Cells.Item(4, 19).Value = '04/08/2023 05:42'
$ws.Cells.Item(4, 20).Value = 5.11
$ws.Cells.Item(4, 21).Value = '05/08/2023 17:25'
$ws.Cells.Item(4, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/jedinstvo-ub-rfk-novi-sad/EFrxoV3l/'

# Row 7
$ws.Cells.Item(7, 6).Value = 'Mladost GAT'
$ws.Cells.Item(7, 7).Value = 1
$ws.Cells.Item(7, 8).Value = 'Radnicki S. Mitrovica'
$ws.Cells.Item(7, 9).Value = 1
$ws.Cells.Item(7, 10).Value = 1.75
$ws.Cells.Item(7, 11).Value = '05/08/2023 22:30'
$ws.Cells.Item(7, 12).Value = 1.65
$ws.Cells.Item(7, 13).Value = '06/08/2023 16:38'
$ws.Cells.Item(7, 14).Value = 3.27
$ws.Cells.Item(7, 15).Value = '05/08/2023 22:30'
$ws.Cells.Item(7, 16).Value = 3.31
$ws.Cells.Item(7, 17).Value = '06/08/2023 16:39'
$ws.Cells.Item(7, 18).Value = 4.28
$ws.Cells.Item(7, 19).Value = '05/08/2023 22:30'
$ws.Cells.Item(7, 20).Value = 5.05
$ws.Cells.Item(7, 21).Value = '06/08/2023 16:39'
$ws.Cells.Item(7, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/mladost-gat-radnicki-s-mitrovica/MH3chDBE/'

# Row 8
$ws.Cells.Item(8, 6).Value = 'Radnicki Beograd'
$ws.Cells.Item(8, 7).Value = 2
$ws.Cells.Item(8, 8).Value = 'Dubocica'
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 1.87
$ws.Cells.Item(8, 11).Value = '05/08/2023 22:30'
$ws.Cells.Item(8, 12).Value = 2.06
$ws.Cells.Item(8, 13).Value = '06/08/2023 16:38'
$ws.Cells.Item(8, 14).Value = 3.12
$ws.Cells.Item(8, 15).Value = '05/08/2023 22:30'
$ws.Cells.Item(8, 16).Value = 3.34
$ws.Cells.Item(8, 17).Value = '06/08/2023 16:36'
$ws.Cells.Item(8, 18).Value = 3.89
$ws.Cells.Item(8, 19).Value = '05/08/2023 22:30'
$ws.Cells.Item(8, 20).Value = 3.11
$ws.Cells.Item(8, 21).Value = '06/08/2023 16:38'
$ws.Cells.Item(8, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/radnicki-beograd-dubocica/fRaZqkJf/'

# Row 15
$ws.Cells.Item(15, 6).Value = 'RFK Novi Sad'
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(15, 8).Value = 'Radnicki Beograd'
$ws.Cells.Item(15, 9).Value = 1
$ws.Cells.Item(15, 10).Value = 1.79
$ws.Cells.Item(15, 11).Value = '12/08/2023 05:43'
$ws.Cells.Item(15, 12).Value = 2.52
$ws.Cells.Item(15, 13).Value = '13/08/2023 16:46'
$ws.Cells.Item(15, 14).Value = 3.13
$ws.Cells.Item(15, 15).Value = '12/08/2023 05:43'
$ws.Cells.Item(15, 16).Value = 3.21
$ws.Cells.Item(15, 17).Value = '13/08/2023 17:01'
$ws.Cells.Item(15, 18).Value = 3.68
$ws.Cells.Item(15, 19).Value = '12/08/2023 05:43'
$ws.Cells.Item(15, 20).Value = 2.52
$ws.Cells.Item(15, 21).Value = '13/08/2023 16:46'
$ws.Cells.Item(15, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/rfk-novi-sad-radnicki-beograd/SvbsrTm7/'

# Row 16
$ws.Cells.Item(16, 6).Value = 'OFK Beograd'
$ws.Cells.Item(16, 7).Value = 4
$ws.Cells.Item(16, 8).Value = 'Metalac'
$ws.Cells.Item(16, 9).Value = 1
$ws.Cells.Item(16, 10).Value = 1.87
$ws.Cells.Item(16, 11).Value = '12/08/2023 23:08'
$ws.Cells.Item(16, 12).Value = 1.83
$ws.Cells.Item(16, 13).Value = '13/08/2023 17:15'
$ws.Cells.Item(16, 14).Value = 3.15
$ws.Cells.Item(16, 15).Value = '12/08/2023 23:08'
$ws.Cells.Item(16, 16).Value = 3.62
$ws.Cells.Item(16, 17).Value = '13/08/2023 17:15'
$ws.Cells.Item(16, 18).Value = 3.85
$ws.Cells.Item(16, 19).Value = '12/08/2023 23:08'
$ws.Cells.Item(16, 20).Value = 3.52
$ws.Cells.Item(16, 21).Value = '13/08/2023 17:15'
$ws.Cells.Item(16, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/ofk-beograd-metalac/rmtU5WuE/'

# Row 18
$ws.Cells.Item(18, 6).Value = 'Macva'
$ws.Cells.Item(18, 7).Value = 3
$ws.Cells.Item(18, 8).Value = 'Radnicki S. Mitrovica'
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 2.14
$ws.Cells.Item(18, 11).Value = '18/08/2023 07:12'
$ws.Cells.Item(18, 12).Value = 2.14
$ws.Cells.Item(18, 13).Value = '18/08/2023 07:12'
$ws.Cells.Item(18, 14).Value = 2.75
$ws.Cells.Item(18, 15).Value = '18/08/2023 07:12'
$ws.Cells.Item(18, 16).Value = 2.83
$ws.Cells.Item(18, 17).Value = '19/08/2023 17:02'
$ws.Cells.Item(18, 18).Value = 3.17
$ws.Cells.Item(18, 19).Value = '18/08/2023 07:12'
$ws.Cells.Item(18, 20).Value = 3.17
$ws.Cells.Item(18, 21).Value = '18/08/2023 07:12'
$ws.Cells.Item(18, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/macva-sabac-radnicki-s-mitrovica/YVwCc8mf/'

# Row 19
$ws.Cells.Item(19, 6).Value = 'Sloboda'
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 'Kolubara'
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 2.58
$ws.Cells.Item(19, 11).Value = '19/08/2023 09:21'
$ws.Cells.Item(19, 12).Value = 2.84
$ws.Cells.Item(19, 13).Value = '19/08/2023 18:40'
$ws.Cells.Item(19, 14).Value = 2.84
$ws.Cells.Item(19, 15).Value = '19/08/2023 09:21'
$ws.Cells.Item(19, 16).Value = 2.84
$ws.Cells.Item(19, 17).Value = '19/08/2023 18:40'
$ws.Cells.Item(19, 18).Value = 2.68
$ws.Cells.Item(19, 19).Value = '19/08/2023 09:21'
$ws.Cells.Item(19, 20).Value = 2.41
$ws.Cells.Item(19, 21).Value = '19/08/2023 18:40'
$ws.Cells.Item(19, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/sloboda-kolubara/bw9EL1AP/'

# Row 20
$ws.Cells.Item(20, 6).Value = 'Metalac'
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 'Tekstilac Odzaci'
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 1.81
$ws.Cells.Item(20, 11).Value = '19/08/2023 09:21'
$ws.Cells.Item(20, 12).Value = 1.81
$ws.Cells.Item(20, 13).Value = '19/08/2023 09:21'
$ws.Cells.Item(20, 14).Value = 3.12
$ws.Cells.Item(20, 15).Value = '19/08/2023 09:21'
$ws.Cells.Item(20, 16).Value = 3.14
$ws.Cells.Item(20, 17).Value = '19/08/2023 17:06'
$ws.Cells.Item(20, 18).Value = 4.18
$ws.Cells.Item(20, 19).Value = '19/08/2023 09:21'
$ws.Cells.Item(20, 20).Value = 4.18
$ws.Cells.Item(20, 21).Value = '19/08/2023 09:21'
$ws.Cells.Item(20, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/metalac-tekstilac-odzaci/ShwGdS20/'

# Row 33
$ws.Cells.Item(33, 6).Value = 'Jedinstvo U.'
$ws.Cells.Item(33, 7).Value = 3
$ws.Cells.Item(33, 8).Value = 'Kolubara'
$ws.Cells.Item(33, 9).Value = 2
$ws.Cells.Item(33, 10).Value = 2.34
$ws.Cells.Item(33, 11).Value = '02/09/2023 15:12'
$ws.Cells.Item(33, 12).Value = 2.12
$ws.Cells.Item(33, 13).Value = '02/09/2023 16:58'
$ws.Cells.Item(33, 14).Value = 2.84
$ws.Cells.Item(33, 15).Value = '02/09/2023 15:12'
$ws.Cells.Item(33, 16).Value = 2.84
$ws.Cells.Item(33, 17).Value = '02/09/2023 16:57'
$ws.Cells.Item(33, 18).Value = 2.99
$ws.Cells.Item(33, 19).Value = '02/09/2023 15:12'
$ws.Cells.Item(33, 20).Value = 3.54
$ws.Cells.Item(33, 21).Value = '02/09/2023 16:58'
$ws.Cells.Item(33, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/jedinstvo-ub-kolubara/8CfpGuYn/'

# Row 34
$ws.Cells.Item(34, 6).Value = 'Vrsac'
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 'FK Indjija'
$ws.Cells.Item(34, 9).Value = 1
$ws.Cells.Item(34, 10).Value = 2.27
$ws.Cells.Item(34, 11).Value = '01/09/2023 05:13'
$ws.Cells.Item(34, 12).Value = 2.71
$ws.Cells.Item(34, 13).Value = '02/09/2023 16:58'
$ws.Cells.Item(34, 14).Value = 2.81
$ws.Cells.Item(34, 15).Value = '01/09/2023 05:13'
$ws.Cells.Item(34, 16).Value = 2.84
$ws.Cells.Item(34, 17).Value = '02/09/2023 16:03'
$ws.Cells.Item(34, 18).Value = 2.86
$ws.Cells.Item(34, 19).Value = '01/09/2023 05:13'
$ws.Cells.Item(34, 20).Value = 2.61
$ws.Cells.Item(34, 21).Value = '02/09/2023 16:58'
$ws.Cells.Item(34, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/vrsac-indjija/zNYfcagB/'

# Row 45
$ws.Cells.Item(45, 6).Value = 'OFK Beograd'
$ws.Cells.Item(45, 7).Value = 4
$ws.Cells.Item(45, 8).Value = 'RFK Novi Sad'
$ws.Cells.Item(45, 9).Value = 1
$ws.Cells.Item(45, 10).Value = 1.22
$ws.Cells.Item(45, 11).Value = '10/09/2023 14:11'
$ws.Cells.Item(45, 12).Value = 1.1
$ws.Cells.Item(45, 13).Value = '10/09/2023 16:27'
$ws.Cells.Item(45, 14).Value = 5.09
$ws.Cells.Item(45, 15).Value = '10/09/2023 14:11'
$ws.Cells.Item(45, 16).Value = 7.91
$ws.Cells.Item(45, 17).Value = '10/09/2023 16:27'
$ws.Cells.Item(45, 18).Value = 10.53
$ws.Cells.Item(45, 19).Value = '10/09/2023 14:11'
$ws.Cells.Item(45, 20).Value = 20.12
$ws.Cells.Item(45, 21).Value = '10/09/2023 16:27'
$ws.Cells.Item(45, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/ofk-beograd-rfk-novi-sad/lQbUBb84/'

# Row 46
$ws.Cells.Item(46, 6).Value = 'Mladost GAT'
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 'Sloboda'
$ws.Cells.Item(46, 9).Value = 1
$ws.Cells.Item(46, 10).Value = 1.63
$ws.Cells.Item(46, 11).Value = '10/09/2023 14:12'
$ws.Cells.Item(46, 12).Value = 1.63
$ws.Cells.Item(46, 13).Value = '10/09/2023 14:12'
$ws.Cells.Item(46, 14).Value = 3.34
$ws.Cells.Item(46, 15).Value = '10/09/2023 14:12'
$ws.Cells.Item(46, 16).Value = 3.36
$ws.Cells.Item(46, 17).Value = '10/09/2023 14:34'
$ws.Cells.Item(46, 18).Value = 5.04
$ws.Cells.Item(46, 19).Value = '10/09/2023 14:12'
$ws.Cells.Item(46, 20).Value = 5.04
$ws.Cells.Item(46, 21).Value = '10/09/2023 14:12'
$ws.Cells.Item(46, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/mladost-gat-sloboda/I5jo8G7T/'

# Row 47
$ws.Cells.Item(47, 6).Value = 'Tekstilac Odzaci'
$ws.Cells.Item(47, 7).Value = 1
$ws.Cells.Item(47, 8).Value = 'Dubocica'
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 2.23
$ws.Cells.Item(47, 11).Value = '10/09/2023 14:10'
$ws.Cells.Item(47, 12).Value = 1.85
$ws.Cells.Item(47, 13).Value = '10/09/2023 16:21'
$ws.Cells.Item(47, 14).Value = 3
$ws.Cells.Item(47, 15).Value = '10/09/2023 14:10'
$ws.Cells.Item(47, 16).Value = 3.11
$ws.Cells.Item(47, 17).Value = '10/09/2023 16:21'
$ws.Cells.Item(47, 18).Value = 3.03
$ws.Cells.Item(47, 19).Value = '10/09/2023 14:10'
$ws.Cells.Item(47, 20).Value = 4.07
$ws.Cells.Item(47, 21).Value = '10/09/2023 16:21'
$ws.Cells.Item(47, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/tekstilac-odzaci-dubocica/WllZAINA/'

# Row 58
$ws.Cells.Item(58, 6).Value = 'Macva'
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 8).Value = 'Sloboda'
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 10).Value = 1.93
$ws.Cells.Item(58, 11).Value = '21/09/2023 06:13'
$ws.Cells.Item(58, 12).Value = 1.74
$ws.Cells.Item(58, 13).Value = '22/09/2023 16:13'
$ws.Cells.Item(58, 14).Value = 2.79
$ws.Cells.Item(58, 15).Value = '21/09/2023 06:13'
$ws.Cells.Item(58, 16).Value = 3.08
$ws.Cells.Item(58, 17).Value = '22/09/2023 17:04'
$ws.Cells.Item(58, 18).Value = 3.65
$ws.Cells.Item(58, 19).Value = '21/09/2023 06:13'
$ws.Cells.Item(58, 20).Value = 4.76
$ws.Cells.Item(58, 21).Value = '22/09/2023 16:13'
$ws.Cells.Item(58, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/macva-sabac-sloboda/pj1uM3m9/'

# Row 59
$ws.Cells.Item(59, 6).Value = 'Kolubara'
$ws.Cells.Item(59, 7).Value = 4
$ws.Cells.Item(59, 8).Value = 'RFK Novi Sad'
$ws.Cells.Item(59, 9).Value = 1
$ws.Cells.Item(59, 10).Value = 1.36
$ws.Cells.Item(59, 11).Value = '21/09/2023 06:13'
$ws.Cells.Item(59, 12).Value = 1.38
$ws.Cells.Item(59, 13).Value = '22/09/2023 18:35'
$ws.Cells.Item(59, 14).Value = 4.04
$ws.Cells.Item(59, 15).Value = '21/09/2023 06:13'
$ws.Cells.Item(59, 16).Value = 4.2
$ws.Cells.Item(59, 17).Value = '22/09/2023 18:35'
$ws.Cells.Item(59, 18).Value = 5.95
$ws.Cells.Item(59, 19).Value = '21/09/2023 06:13'
$ws.Cells.Item(59, 20).Value = 7.07
$ws.Cells.Item(59, 21).Value = '22/09/2023 18:35'
$ws.Cells.Item(59, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/kolubara-rfk-novi-sad/rVncDwI5/'

# Row 73
$ws.Cells.Item(73, 6).Value = 'Tekstilac Odzaci'
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 'Kolubara'
$ws.Cells.Item(73, 9).Value = 1
$ws.Cells.Item(73, 10).Value = 2.19
$ws.Cells.Item(73, 11).Value = '07/10/2023 02:14'
$ws.Cells.Item(73, 12).Value = 1.98
$ws.Cells.Item(73, 13).Value = '08/10/2023 14:46'
$ws.Cells.Item(73, 14).Value = 2.76
$ws.Cells.Item(73, 15).Value = '07/10/2023 02:14'
$ws.Cells.Item(73, 16).Value = 2.93
$ws.Cells.Item(73, 17).Value = '08/10/2023 14:45'
$ws.Cells.Item(73, 18).Value = 3.04
$ws.Cells.Item(73, 19).Value = '07/10/2023 02:14'
$ws.Cells.Item(73, 20).Value = 3.85
$ws.Cells.Item(73, 21).Value = '08/10/2023 14:46'
$ws.Cells.Item(73, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/tekstilac-odzaci-kolubara/YJuAAy2O/'

# Row 74
$ws.Cells.Item(74, 6).Value = 'OFK Beograd'
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 'Radnicki S. Mitrovica'
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 1.38
$ws.Cells.Item(74, 11).Value = '06/10/2023 02:12'
$ws.Cells.Item(74, 12).Value = 1.35
$ws.Cells.Item(74, 13).Value = '08/10/2023 14:36'
$ws.Cells.Item(74, 14).Value = 3.92
$ws.Cells.Item(74, 15).Value = '06/10/2023 02:12'
$ws.Cells.Item(74, 16).Value = 4.31
$ws.Cells.Item(74, 17).Value = '08/10/2023 14:36'
$ws.Cells.Item(74, 18).Value = 5.88
$ws.Cells.Item(74, 19).Value = '06/10/2023 02:12'
$ws.Cells.Item(74, 20).Value = 7.42
$ws.Cells.Item(74, 21).Value = '08/10/2023 14:36'
$ws.Cells.Item(74, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/ofk-beograd-radnicki-s-mitrovica/pWGXTiom/'

# Row 75
$ws.Cells.Item(75, 6).Value = 'Radnicki Beograd'
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 'Macva'
$ws.Cells.Item(75, 9).Value = 1
$ws.Cells.Item(75, 10).Value = 2.56
$ws.Cells.Item(75, 11).Value = '07/10/2023 02:14'
$ws.Cells.Item(75, 12).Value = 2.87
$ws.Cells.Item(75, 13).Value = '08/10/2023 14:41'
$ws.Cells.Item(75, 14).Value = 2.67
$ws.Cells.Item(75, 15).Value = '07/10/2023 02:14'
$ws.Cells.Item(75, 16).Value = 2.75
$ws.Cells.Item(75, 17).Value = '08/10/2023 14:41'
$ws.Cells.Item(75, 18).Value = 2.63
$ws.Cells.Item(75, 19).Value = '07/10/2023 02:14'
$ws.Cells.Item(75, 20).Value = 2.55
$ws.Cells.Item(75, 21).Value = '08/10/2023 14:41'
$ws.Cells.Item(75, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/radnicki-beograd-macva-sabac/YimFBAVI/'

# Row 84
$ws.Cells.Item(84, 6).Value = 'Mladost GAT'
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 'Graficar Beograd'
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 1.97
$ws.Cells.Item(84, 11).Value = '15/10/2023 02:12'
$ws.Cells.Item(84, 12).Value = 1.88
$ws.Cells.Item(84, 13).Value = '16/10/2023 14:54'
$ws.Cells.Item(84, 14).Value = 3.07
$ws.Cells.Item(84, 15).Value = '15/10/2023 02:12'
$ws.Cells.Item(84, 16).Value = 3.51
$ws.Cells.Item(84, 17).Value = '16/10/2023 14:54'
$ws.Cells.Item(84, 18).Value = 3.17
$ws.Cells.Item(84, 19).Value = '15/10/2023 02:12'
$ws.Cells.Item(84, 20).Value = 3.45
$ws.Cells.Item(84, 21).Value = '16/10/2023 14:54'
$ws.Cells.Item(84, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/mladost-gat-graficar-beograd/zL9pRkV5/'

# Row 86
$ws.Cells.Item(86, 6).Value = 'Vrsac'
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = 'RFK Novi Sad'
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 1.51
$ws.Cells.Item(86, 11).Value = '15/10/2023 02:12'
$ws.Cells.Item(86, 12).Value = 1.44
$ws.Cells.Item(86, 13).Value = '16/10/2023 14:07'
$ws.Cells.Item(86, 14).Value = 3.5
$ws.Cells.Item(86, 15).Value = '15/10/2023 02:12'
$ws.Cells.Item(86, 16).Value = 3.72
$ws.Cells.Item(86, 17).Value = '16/10/2023 14:07'
$ws.Cells.Item(86, 18).Value = 5.01
$ws.Cells.Item(86, 19).Value = '15/10/2023 02:12'
$ws.Cells.Item(86, 20).Value = 7.01
$ws.Cells.Item(86, 21).Value = '16/10/2023 14:07'
$ws.Cells.Item(86, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/vrsac-rfk-novi-sad/EB8lQ9pC/'

# Row 87
$ws.Cells.Item(87, 6).Value = 'Radnicki Beograd'
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 'Metalac'
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 2.72
$ws.Cells.Item(87, 11).Value = '20/10/2023 01:13'
$ws.Cells.Item(87, 12).Value = 3.96
$ws.Cells.Item(87, 13).Value = '21/10/2023 13:58'
$ws.Cells.Item(87, 14).Value = 2.69
$ws.Cells.Item(87, 15).Value = '20/10/2023 01:13'
$ws.Cells.Item(87, 16).Value = 2.62
$ws.Cells.Item(87, 17).Value = '21/10/2023 13:46'
$ws.Cells.Item(87, 18).Value = 2.46
$ws.Cells.Item(87, 19).Value = '20/10/2023 01:13'
$ws.Cells.Item(87, 20).Value = 2.13
$ws.Cells.Item(87, 21).Value = '21/10/2023 13:58'
$ws.Cells.Item(87, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/radnicki-beograd-metalac/61g4sBFn/'

# Row 88
$ws.Cells.Item(88, 6).Value = 'Jedinstvo U.'
$ws.Cells.Item(88, 7).Value = 1
$ws.Cells.Item(88, 8).Value = 'FK Indjija'
$ws.Cells.Item(88, 9).Value = 1
$ws.Cells.Item(88, 10).Value = 2.11
$ws.Cells.Item(88, 11).Value = '20/10/2023 07:12'
$ws.Cells.Item(88, 12).Value = 2.37
$ws.Cells.Item(88, 13).Value = '21/10/2023 13:45'
$ws.Cells.Item(88, 14).Value = 2.81
$ws.Cells.Item(88, 15).Value = '20/10/2023 07:12'
$ws.Cells.Item(88, 16).Value = 2.79
$ws.Cells.Item(88, 17).Value = '21/10/2023 13:45'
$ws.Cells.Item(88, 18).Value = 3.16
$ws.Cells.Item(88, 19).Value = '20/10/2023 07:12'
$ws.Cells.Item(88, 20).Value = 3.09
$ws.Cells.Item(88, 21).Value = '21/10/2023 13:45'
$ws.Cells.Item(88, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/jedinstvo-ub-indjija/Ox4ari0t/'

# Row 89
$ws.Cells.Item(89, 6).Value = 'OFK Beograd'
$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(89, 8).Value = 'Mladost GAT'
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 1.43
$ws.Cells.Item(89, 11).Value = '20/10/2023 01:13'
$ws.Cells.Item(89, 12).Value = 1.29
$ws.Cells.Item(89, 13).Value = '21/10/2023 13:58'
$ws.Cells.Item(89, 14).Value = 3.77
$ws.Cells.Item(89, 15).Value = '20/10/2023 01:13'
$ws.Cells.Item(89, 16).Value = 4.64
$ws.Cells.Item(89, 17).Value = '21/10/2023 13:58'
$ws.Cells.Item(89, 18).Value = 5.37
$ws.Cells.Item(89, 19).Value = '20/10/2023 01:13'
$ws.Cells.Item(89, 20).Value = 8.99
$ws.Cells.Item(89, 21).Value = '21/10/2023 13:58'
$ws.Cells.Item(89, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/ofk-beograd-mladost-gat/zsoLwTEB/'

# Row 108
$ws.Cells.Item(108, 6).Value = 'Radnicki S. Mitrovica'
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 'Smederevo'
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 1.92
$ws.Cells.Item(108, 11).Value = '06/11/2023 01:12'
$ws.Cells.Item(108, 12).Value = 1.78
$ws.Cells.Item(108, 13).Value = '06/11/2023 12:25'
$ws.Cells.Item(108, 14).Value = 2.95
$ws.Cells.Item(108, 15).Value = '06/11/2023 01:12'
$ws.Cells.Item(108, 16).Value = 3.21
$ws.Cells.Item(108, 17).Value = '06/11/2023 12:31'
$ws.Cells.Item(108, 18).Value = 3.83
$ws.Cells.Item(108, 19).Value = '06/11/2023 01:12'
$ws.Cells.Item(108, 20).Value = 4.23
$ws.Cells.Item(108, 21).Value = '06/11/2023 12:25'
$ws.Cells.Item(108, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/radnicki-s-mitrovica-smederevo/CtiO77Ti/'

# Row 109
$ws.Cells.Item(109, 6).Value = 'Radnicki Beograd'
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = 'FK Indjija'
$ws.Cells.Item(109, 9).Value = 1
$ws.Cells.Item(109, 10).Value = 2.65
$ws.Cells.Item(109, 11).Value = '06/11/2023 01:12'
$ws.Cells.Item(109, 12).Value = 2.75
$ws.Cells.Item(109, 13).Value = '06/11/2023 12:57'
$ws.Cells.Item(109, 14).Value = 2.86
$ws.Cells.Item(109, 15).Value = '06/11/2023 01:12'
$ws.Cells.Item(109, 16).Value = 2.72
$ws.Cells.Item(109, 17).Value = '06/11/2023 12:58'
$ws.Cells.Item(109, 18).Value = 2.54
$ws.Cells.Item(109, 19).Value = '06/11/2023 01:12'
$ws.Cells.Item(109, 20).Value = 2.68
$ws.Cells.Item(109, 21).Value = '06/11/2023 12:58'
$ws.Cells.Item(109, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/radnicki-beograd-indjija/SzAbinqB/'

# ---------------------------------------------------------------------------
# Part 2: append 3 brand-new match rows (114-116), copying number formats
# from the last pre-existing data row (113) for the styled columns A and E.
# ---------------------------------------------------------------------------

# Row 114
$ws.Cells.Item(114, 1).Value = 113
$ws.Cells.Item(114, 2).Value = 'serbia'
$ws.Cells.Item(114, 3).Value = 'prva-liga'
$ws.Cells.Item(114, 4).Value = '2023-2024'
$ws.Cells.Item(114, 5).Value = 45242.54166666666
$ws.Cells.Item(114, 6).Value = 'OFK Beograd'
$ws.Cells.Item(114, 7).Value = 1
$ws.Cells.Item(114, 8).Value = 'Sloboda'
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 1.26
$ws.Cells.Item(114, 11).Value = '12/11/2023 02:12'
$ws.Cells.Item(114, 12).Value = 1.32
$ws.Cells.Item(114, 13).Value = '12/11/2023 12:21'
$ws.Cells.Item(114, 14).Value = 4.79
$ws.Cells.Item(114, 15).Value = '12/11/2023 02:12'
$ws.Cells.Item(114, 16).Value = 4.55
$ws.Cells.Item(114, 17).Value = '12/11/2023 12:21'
$ws.Cells.Item(114, 18).Value = 8.58
$ws.Cells.Item(114, 19).Value = '12/11/2023 02:12'
$ws.Cells.Item(114, 20).Value = 8.03
$ws.Cells.Item(114, 21).Value = '12/11/2023 12:21'
$ws.Cells.Item(114, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/ofk-beograd-sloboda/E5mu4PSG/'
$ws.Cells.Item(113, 1).Copy() | Out-Null
$ws.Cells.Item(114, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(113, 5).Copy() | Out-Null
$ws.Cells.Item(114, 5).PasteSpecial(-4122) | Out-Null

# Row 115
$ws.Cells.Item(115, 1).Value = 114
$ws.Cells.Item(115, 2).Value = 'serbia'
$ws.Cells.Item(115, 3).Value = 'prva-liga'
$ws.Cells.Item(115, 4).Value = '2023-2024'
$ws.Cells.Item(115, 5).Value = 45242.54166666666
$ws.Cells.Item(115, 6).Value = 'Graficar Beograd'
$ws.Cells.Item(115, 7).Value = 1
$ws.Cells.Item(115, 8).Value = 'Metalac'
$ws.Cells.Item(115, 9).Value = 1
$ws.Cells.Item(115, 10).Value = 1.75
$ws.Cells.Item(115, 11).Value = '12/11/2023 02:12'
$ws.Cells.Item(115, 12).Value = 1.78
$ws.Cells.Item(115, 13).Value = '12/11/2023 12:48'
$ws.Cells.Item(115, 14).Value = 3.23
$ws.Cells.Item(115, 15).Value = '12/11/2023 02:12'
$ws.Cells.Item(115, 16).Value = 3.26
$ws.Cells.Item(115, 17).Value = '12/11/2023 12:48'
$ws.Cells.Item(115, 18).Value = 4.14
$ws.Cells.Item(115, 19).Value = '12/11/2023 02:12'
$ws.Cells.Item(115, 20).Value = 4.21
$ws.Cells.Item(115, 21).Value = '12/11/2023 12:48'
$ws.Cells.Item(115, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/graficar-beograd-metalac/f319nccq/'
$ws.Cells.Item(113, 1).Copy() | Out-Null
$ws.Cells.Item(115, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(113, 5).Copy() | Out-Null
$ws.Cells.Item(115, 5).PasteSpecial(-4122) | Out-Null

# Row 116
$ws.Cells.Item(116, 1).Value = 115
$ws.Cells.Item(116, 2).Value = 'serbia'
$ws.Cells.Item(116, 3).Value = 'prva-liga'
$ws.Cells.Item(116, 4).Value = '2023-2024'
$ws.Cells.Item(116, 5).Value = 45242.54166666666
$ws.Cells.Item(116, 6).Value = 'RFK Novi Sad'
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 'Jedinstvo U.'
$ws.Cells.Item(116, 9).Value = 2
$ws.Cells.Item(116, 10).Value = 3.06
$ws.Cells.Item(116, 11).Value = '12/11/2023 02:12'
$ws.Cells.Item(116, 12).Value = 3.08
$ws.Cells.Item(116, 13).Value = '12/11/2023 12:20'
$ws.Cells.Item(116, 14).Value = 2.92
$ws.Cells.Item(116, 15).Value = '12/11/2023 02:12'
$ws.Cells.Item(116, 16).Value = 3.17
$ws.Cells.Item(116, 17).Value = '12/11/2023 12:20'
$ws.Cells.Item(116, 18).Value = 2.22
$ws.Cells.Item(116, 19).Value = '12/11/2023 02:12'
$ws.Cells.Item(116, 20).Value = 2.15
$ws.Cells.Item(116, 21).Value = '12/11/2023 12:20'
$ws.Cells.Item(116, 22).Value = 'https://www.betexplorer.com/football/serbia/prva-liga/rfk-novi-sad-jedinstvo-ub/Yg0DoHCk/'
$ws.Cells.Item(113, 1).Copy() | Out-Null
$ws.Cells.Item(116, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(113, 5).Copy() | Out-Null
$ws.Cells.Item(116, 5).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

Write-Host "Edit applied: reordered 27 same-date rows and appended rows 114-116."
